$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value2 = 1.68
$ws.Cells.Item(2, 7).Value2 = 1.81
$ws.Cells.Item(2, 8).Value2 = 5.2
$ws.Cells.Item(2, 9).Value2 = 6.2
$ws.Cells.Item(2, 10).Value2 = 3.75
$ws.Cells.Item(2, 11).Value2 = 4.4
$ws.Cells.Item(2, 12).Value2 = 1.37
$ws.Cells.Item(2, 14).Value2 = 3.7
$ws.Cells.Item(2, 15).Value2 = 1.3
$ws.Cells.Item(2, 17).Value2 = 1.89
$ws.Cells.Item(2, 18).Value2 = 1.38
$ws.Cells.Item(2, 21).Value2 = 1.96
$ws.Cells.Item(2, 23).Value2 = 2.22
$ws.Cells.Item(2, 25).Value2 = 1000
$ws.Cells.Item(2, 28).Value2 = 29
$ws.Cells.Item(2, 29).Value2 = 14
$ws.Cells.Item(2, 32).Value2 = 40
$ws.Cells.Item(2, 33).Value2 = 40
$ws.Cells.Item(2, 40).Value2 = 85
$ws.Cells.Item(3, 8).Value2 = 4.4
$ws.Cells.Item(3, 10).Value2 = 3.85
$ws.Cells.Item(3, 18).Value2 = 1.46
$ws.Cells.Item(3, 20).Value2 = 1.68
$ws.Cells.Item(3, 21).Value2 = 2.16
$ws.Cells.Item(3, 24).Value2 = 24
$ws.Cells.Item(3, 26).Value2 = 1000
$ws.Cells.Item(3, 29).Value2 = 10.5
$ws.Cells.Item(3, 30).Value2 = 24
$ws.Cells.Item(3, 31).Value2 = 220
$ws.Cells.Item(3, 32).Value2 = 13
$ws.Cells.Item(3, 35).Value2 = 170
$ws.Cells.Item(3, 36).Value2 = 22
$ws.Cells.Item(3, 38).Value2 = 36
$ws.Cells.Item(3, 41).Value2 = 300
$ws.Cells.Item(4, 10).Value2 = 2.44
$ws.Cells.Item(4, 14).Value2 = 1.32
$ws.Cells.Item(4, 16).Value2 = 1.32
$ws.Cells.Item(4, 18).Value2 = 1.21
$ws.Cells.Item(4, 20).Value2 = 1.03
$ws.Cells.Item(4, 21).Value2 = 1.03
$ws.Cells.Item(5, 6).Value2 = 1.41
$ws.Cells.Item(5, 7).Value2 = 1.43
$ws.Cells.Item(5, 8).Value2 = 9.4
$ws.Cells.Item(5, 10).Value2 = 5.1
$ws.Cells.Item(5, 11).Value2 = 5.3
$ws.Cells.Item(5, 12).Value2 = 1.37
$ws.Cells.Item(5, 13).Value2 = 1.06
$ws.Cells.Item(5, 15).Value2 = 1.29
$ws.Cells.Item(5, 17).Value2 = 1.87
$ws.Cells.Item(5, 20).Value2 = 2.16
$ws.Cells.Item(5, 21).Value2 = 1.79
$ws.Cells.Item(5, 23).Value2 = 3.3
$ws.Cells.Item(5, 24).Value2 = 16.5
$ws.Cells.Item(5, 25).Value2 = 29
$ws.Cells.Item(5, 26).Value2 = 85
$ws.Cells.Item(5, 27).Value2 = 380
$ws.Cells.Item(5, 30).Value2 = 36
$ws.Cells.Item(5, 31).Value2 = 180
$ws.Cells.Item(5, 33).Value2 = 10
$ws.Cells.Item(5, 35).Value2 = 170
$ws.Cells.Item(5, 36).Value2 = 11.5
$ws.Cells.Item(5, 39).Value2 = 200
$ws.Cells.Item(5, 40).Value2 = 7
$ws.Cells.Item(5, 41).Value2 = 270
$ws.Cells.Item(6, 6).Value2 = 4.3
$ws.Cells.Item(6, 7).Value2 = 4.6
$ws.Cells.Item(6, 9).Value2 = 2.14
$ws.Cells.Item(6, 14).Value2 = 2.96
$ws.Cells.Item(6, 15).Value2 = 1.47
$ws.Cells.Item(6, 16).Value2 = 1.65
$ws.Cells.Item(6, 17).Value2 = 2.4
$ws.Cells.Item(6, 19).Value2 = 4.7
$ws.Cells.Item(6, 21).Value2 = 1.86
$ws.Cells.Item(6, 22).Value2 = 1.88
$ws.Cells.Item(6, 27).Value2 = 30
$ws.Cells.Item(6, 33).Value2 = 20
$ws.Cells.Item(6, 34).Value2 = 23
$ws.Cells.Item(6, 35).Value2 = 200
$ws.Cells.Item(6, 36).Value2 = 120
$ws.Cells.Item(6, 37).Value2 = 160
$ws.Cells.Item(6, 38).Value2 = 95
$ws.Cells.Item(6, 39).Value2 = 580
$ws.Cells.Item(6, 40).Value2 = 260
$ws.Cells.Item(6, 41).Value2 = 24
$ws.Cells.Item(7, 6).Value2 = 2.26
$ws.Cells.Item(7, 7).Value2 = 2.3
$ws.Cells.Item(7, 9).Value2 = 3.85
$ws.Cells.Item(7, 10).Value2 = 3.3
$ws.Cells.Item(7, 11).Value2 = 3.4
$ws.Cells.Item(7, 12).Value2 = 1.48
$ws.Cells.Item(7, 13).Value2 = 1.09
$ws.Cells.Item(7, 14).Value2 = 3.3
$ws.Cells.Item(7, 15).Value2 = 1.4
$ws.Cells.Item(7, 16).Value2 = 1.76
$ws.Cells.Item(7, 17).Value2 = 2.18
$ws.Cells.Item(7, 18).Value2 = 1.28
$ws.Cells.Item(7, 19).Value2 = 4.1
$ws.Cells.Item(7, 20).Value2 = 1.89
$ws.Cells.Item(7, 21).Value2 = 2.02
$ws.Cells.Item(7, 22).Value2 = 1.35
$ws.Cells.Item(7, 23).Value2 = 1.76
$ws.Cells.Item(7, 24).Value2 = 13
$ws.Cells.Item(7, 25).Value2 = 13
$ws.Cells.Item(7, 26).Value2 = 29
$ws.Cells.Item(7, 27).Value2 = 80
$ws.Cells.Item(7, 28).Value2 = 9
$ws.Cells.Item(7, 29).Value2 = 7.8
$ws.Cells.Item(7, 30).Value2 = 17.5
$ws.Cells.Item(7, 31).Value2 = 170
$ws.Cells.Item(7, 32).Value2 = 15.5
$ws.Cells.Item(7, 34).Value2 = 22
$ws.Cells.Item(7, 35).Value2 = 1000
$ws.Cells.Item(7, 36).Value2 = 32
$ws.Cells.Item(7, 37).Value2 = 29
$ws.Cells.Item(7, 39).Value2 = 150
$ws.Cells.Item(7, 40).Value2 = 25
$ws.Cells.Item(7, 41).Value2 = 70
$ws.Cells.Item(8, 6).Value2 = 2.48
$ws.Cells.Item(8, 7).Value2 = 2.56
$ws.Cells.Item(8, 8).Value2 = 3.25
$ws.Cells.Item(8, 11).Value2 = 3.3
$ws.Cells.Item(8, 12).Value2 = 1.54
$ws.Cells.Item(8, 13).Value2 = 1.1
$ws.Cells.Item(8, 14).Value2 = 3
$ws.Cells.Item(8, 15).Value2 = 1.46
$ws.Cells.Item(8, 16).Value2 = 1.65
$ws.Cells.Item(8, 19).Value2 = 4.7
$ws.Cells.Item(8, 20).Value2 = 2
$ws.Cells.Item(8, 21).Value2 = 1.92
$ws.Cells.Item(8, 22).Value2 = 1.41
$ws.Cells.Item(8, 23).Value2 = 1.64
$ws.Cells.Item(8, 24).Value2 = 10
$ws.Cells.Item(8, 26).Value2 = 21
$ws.Cells.Item(8, 27).Value2 = 65
$ws.Cells.Item(8, 28).Value2 = 8.6
$ws.Cells.Item(8, 30).Value2 = 14.5
$ws.Cells.Item(8, 31).Value2 = 48
$ws.Cells.Item(8, 32).Value2 = 15
$ws.Cells.Item(8, 33).Value2 = 12.5
$ws.Cells.Item(8, 36).Value2 = 36
$ws.Cells.Item(8, 37).Value2 = 34
$ws.Cells.Item(8, 39).Value2 = 150
$ws.Cells.Item(8, 40).Value2 = 34
$ws.Cells.Item(8, 41).Value2 = 60
$ws.Cells.Item(9, 6).Value2 = 2.54
$ws.Cells.Item(9, 9).Value2 = 3.35
$ws.Cells.Item(9, 10).Value2 = 3.15
$ws.Cells.Item(9, 11).Value2 = 3.25
$ws.Cells.Item(9, 12).Value2 = 1.52
$ws.Cells.Item(9, 13).Value2 = 1.11
$ws.Cells.Item(9, 14).Value2 = 2.84
$ws.Cells.Item(9, 22).Value2 = 1.42
$ws.Cells.Item(9, 25).Value2 = 11
$ws.Cells.Item(9, 26).Value2 = 20
$ws.Cells.Item(9, 27).Value2 = 65
$ws.Cells.Item(9, 28).Value2 = 9.2
$ws.Cells.Item(9, 29).Value2 = 8.6
$ws.Cells.Item(9, 30).Value2 = 16.5
$ws.Cells.Item(9, 31).Value2 = 980
$ws.Cells.Item(9, 33).Value2 = 13.5
$ws.Cells.Item(9, 35).Value2 = 65
$ws.Cells.Item(10, 6).Value2 = 2.34
$ws.Cells.Item(10, 8).Value2 = 3.6
$ws.Cells.Item(10, 9).Value2 = 3.9
$ws.Cells.Item(10, 13).Value2 = 1.13
$ws.Cells.Item(10, 14).Value2 = 2.58
$ws.Cells.Item(10, 16).Value2 = 1.53
$ws.Cells.Item(10, 17).Value2 = 2.7
$ws.Cells.Item(10, 20).Value2 = 2.16
$ws.Cells.Item(10, 21).Value2 = 1.83
$ws.Cells.Item(10, 22).Value2 = 1.34
$ws.Cells.Item(10, 26).Value2 = 29
$ws.Cells.Item(10, 28).Value2 = 7.2
$ws.Cells.Item(10, 30).Value2 = 19.5
$ws.Cells.Item(10, 33).Value2 = 14.5
$ws.Cells.Item(10, 34).Value2 = 28
$ws.Cells.Item(10, 36).Value2 = 42
$ws.Cells.Item(11, 9).Value2 = 3.6
$ws.Cells.Item(11, 14).Value2 = 3.2
$ws.Cells.Item(11, 16).Value2 = 1.76
$ws.Cells.Item(11, 19).Value2 = 4.1
$ws.Cells.Item(11, 25).Value2 = 13.5
$ws.Cells.Item(11, 29).Value2 = 9
$ws.Cells.Item(11, 31).Value2 = 50
$ws.Cells.Item(11, 32).Value2 = 16.5
$ws.Cells.Item(11, 33).Value2 = 13
$ws.Cells.Item(11, 35).Value2 = 60
$ws.Cells.Item(11, 36).Value2 = 980
$ws.Cells.Item(11, 38).Value2 = 55
$ws.Cells.Item(11, 40).Value2 = 29
$ws.Cells.Item(11, 41).Value2 = 1000
$ws.Cells.Item(12, 6).Value2 = 3.1
$ws.Cells.Item(12, 7).Value2 = 3.2
$ws.Cells.Item(12, 8).Value2 = 2.66
$ws.Cells.Item(12, 9).Value2 = 2.76
$ws.Cells.Item(12, 20).Value2 = 2.12
$ws.Cells.Item(12, 21).Value2 = 1.79
$ws.Cells.Item(12, 22).Value2 = 1.57
$ws.Cells.Item(12, 23).Value2 = 1.45
$ws.Cells.Item(12, 25).Value2 = 8.2
$ws.Cells.Item(12, 26).Value2 = 19
$ws.Cells.Item(12, 28).Value2 = 8.8
$ws.Cells.Item(12, 29).Value2 = 7.4
$ws.Cells.Item(12, 31).Value2 = 980
$ws.Cells.Item(12, 32).Value2 = 19
$ws.Cells.Item(12, 36).Value2 = 430
$ws.Cells.Item(12, 38).Value2 = 70
$ws.Cells.Item(12, 41).Value2 = 40
$ws.Cells.Item(14, 7).Value2 = 2.38
$ws.Cells.Item(14, 9).Value2 = 3.9
$ws.Cells.Item(14, 10).Value2 = 3.25
$ws.Cells.Item(14, 11).Value2 = 3.75
$ws.Cells.Item(14, 12).Value2 = 1.48
$ws.Cells.Item(14, 13).Value2 = 1.09
$ws.Cells.Item(14, 14).Value2 = 2.96
$ws.Cells.Item(14, 15).Value2 = 1.41
$ws.Cells.Item(14, 16).Value2 = 1.71
$ws.Cells.Item(14, 17).Value2 = 2.18
$ws.Cells.Item(14, 18).Value2 = 1.25
$ws.Cells.Item(14, 20).Value2 = 1.89
$ws.Cells.Item(14, 21).Value2 = 1.91
$ws.Cells.Item(14, 22).Value2 = 1.36
$ws.Cells.Item(14, 23).Value2 = 1.72
$ws.Cells.Item(14, 24).Value2 = 12.5
$ws.Cells.Item(14, 28).Value2 = 9.2
$ws.Cells.Item(14, 36).Value2 = 980
$ws.Cells.Item(14, 40).Value2 = 30
$ws.Cells.Item(14, 41).Value2 = 1000
